# Auto commit at 2025-09-24 7:47:41.53
# Updates the Metrics sheet's daily figures (B2:B13). The "today" sheet's
# B11:B22 / E11:E22 / F11:F22 formulas reference these cells directly
# (=Metrics!B2 ... =Metrics!B13) so they recalc automatically once the
# Metrics values change. Finally, restore the recorded cell selections on
# both sheets.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 351135.52999999997
$metrics.Range("B3").Value = 283606.41000000003
$metrics.Range("B4").Value = 110911.29
$metrics.Range("B5").Value = 13863
$metrics.Range("B6").Value = 4270386.4099999992
$metrics.Range("B7").Value = 3611133.8899999997
$metrics.Range("B8").Value = 1240276.9700000002
$metrics.Range("B9").Value = 165023
$metrics.Range("B10").Value = 32735710.210999824
$metrics.Range("B11").Value = 19641003.960000005
$metrics.Range("B12").Value = 11521985.859999999
$metrics.Range("B13").Value = 1262650

$metrics.Range("E38").Select()

$today = $wb.Worksheets.Item("today")
$today.Range("E4").Select()
